$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new bold "Play Faust for Free - Dark Literary-Themed
#    Slot Game" paragraph right before the very last paragraph (the
#    one that used to hold the image-generation prompt).
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($n - 1)
$insertionPoint = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Faust for Free - Dark Literary-Themed Slot Game</w:t></w:r></w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# ------------------------------------------------------------------
# 3) Swap the text of the last paragraph (still italic) from the old
#    image-prompt copy to the meta-description copy. Assigning to
#    .Text (instead of Find/Replace) keeps the straight apostrophe
#    untouched by autocorrect.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$lastRange = $lastPara.Range
$lastRange.MoveEnd(1, -1)
$lastRange.Text = "Discover the formula to defeat the devil in Faust, an immersive slot game inspired by Goethe's masterpiece. Play for free and access bonus features."
